$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: shift header cells M1:R1 -> N1:S1, then set new M1 ---
$ws.Range("S1").Value = $ws.Range("R1").Value2
$ws.Range("R1").Value = $ws.Range("Q1").Value2
$ws.Range("Q1").Value = $ws.Range("P1").Value2
$ws.Range("P1").Value = $ws.Range("O1").Value2
$ws.Range("O1").Value = $ws.Range("N1").Value2
$ws.Range("N1").Value = $ws.Range("M1").Value2
$ws.Range("M1").Value = "SIDEBAR_SUBMENU_SUBMENU"

# --- Row 2: shift data cells L2:M2 -> M2:N2, then set new L2 ---
$ws.Range("N2").Value = $ws.Range("M2").Value2
$ws.Range("M2").Value = $ws.Range("L2").Value2
$ws.Range("L2").Value = "Setup Kelengkapan Kepesertaan"

# Match L2's style to K2's (quote-prefixed left/center style) via format-only paste
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths: new column L (12) gets a plain (non-bestfit) width,
#     the rest keep their former widths, shifted one slot to the right ---
$ws.Columns("L").ColumnWidth = 15
$ws.Columns("M:T").AutoFit()

# --- Selection cursor moved to M13 (as recorded by the saved view state) ---
$ws.Range("M13").Select()
